$d = $word.ActiveDocument

# Locate the "VALOR:" label in the certificate body and insert the new
# " PRECIO_MOTO" placeholder right after it (before the existing space
# that precedes the MOTO_MONEDA placeholder), so the line reads:
#   VALOR: PRECIO_MOTO MOTO_MONEDA
$rng = $d.Content
$found = $rng.Find.Execute("VALOR:", $true, $false, $false, $false, $false,
                            $true, 1, $false, "", 0)

if ($found) {
    $insertPoint = $d.Range($rng.End, $rng.End)
    $insertPoint.InsertAfter(" PRECIO_MOTO")
}
